$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.02293332269107
$ws.Range("D2").Value = 1.02674372820799
$ws.Range("E2").Value = 1.047394025025249
$ws.Range("F2").Value = 1.051408384591935
$ws.Range("I2").Value = 1.028621202615219
$ws.Range("J2").Value = 1.028116376318999
$ws.Range("K2").Value = 1.029565482555791
$ws.Range("L2").Value = 1.050156785927615
$ws.Range("M2").Value = 1.054159969073384
$ws.Range("N2").Value = 1.013411189061657
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.024105999484649
$ws.Range("D3").Value = 1.027585047746423
$ws.Range("E3").Value = 1.048776641737936
$ws.Range("F3").Value = 1.052896492310152
$ws.Range("I3").Value = 1.028841044442913
$ws.Range("J3").Value = 1.028926157244195
$ws.Range("K3").Value = 1.03021428950954
$ws.Range("L3").Value = 1.05134975834112
$ws.Range("M3").Value = 1.055458980717611
$ws.Range("N3").Value = 1.013679957185488
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.02486359861419
$ws.Range("D4").Value = 1.028128031750946
$ws.Range("E4").Value = 1.049671326840903
$ws.Range("F4").Value = 1.053859333570421
$ws.Range("I4").Value = 1.028980883488129
$ws.Range("J4").Value = 1.029448458789549
$ws.Range("K4").Value = 1.03063206822088
$ws.Range("L4").Value = 1.052121202791185
$ws.Range("M4").Value = 1.056298957288539
$ws.Range("N4").Value = 1.013853259973024
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.025181808379336
$ws.Range("D5").Value = 1.028355966516706
$ws.Range("E5").Value = 1.050047466104768
$ws.Range("F5").Value = 1.054264101218934
$ws.Range("I5").Value = 1.029039094298536
$ws.Range("J5").Value = 1.029667633518022
$ws.Range("K5").Value = 1.030807214288098
$ws.Range("L5").Value = 1.052445404812768
$ws.Range("M5").Value = 1.056651950700688
$ws.Range("N5").Value = 1.013925971346456
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.025235220568586
$ws.Range("D6").Value = 1.028394218105998
$ws.Range("E6").Value = 1.050110622463105
$ws.Range("F6").Value = 1.054332062984061
$ws.Range("I6").Value = 1.029048834289796
$ws.Range("J6").Value = 1.029704410458766
$ws.Range("K6").Value = 1.030836593487766
$ws.Range("L6").Value = 1.052499833226954
$ws.Range("M6").Value = 1.0567112121776
$ws.Range("N6").Value = 1.013938171409898
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.024867851664886
$ws.Range("D7").Value = 1.028131078744508
$ws.Range("E7").Value = 1.049676352773366
$ws.Range("F7").Value = 1.053864742132482
$ws.Range("I7").Value = 1.028981663572859
$ws.Range("J7").Value = 1.02945138898614
$ws.Range("K7").Value = 1.030634410448003
$ws.Range("L7").Value = 1.052125535234749
$ws.Range("M7").Value = 1.056303674520504
$ws.Range("N7").Value = 1.013854232115594
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.023329884618626
$ws.Range("D8").Value = 1.027028347747758
$ws.Range("E8").Value = 1.047861281115638
$ws.Range("F8").Value = 1.051911313509083
$ws.Range("I8").Value = 1.028695999026311
$ws.Range("J8").Value = 1.028390394553529
$ws.Range("K8").Value = 1.029785173446405
$ws.Range("L8").Value = 1.050560059922188
$ws.Range("M8").Value = 1.054599096696645
$ws.Range("N8").Value = 1.013502146720088
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.020610470499326
$ws.Range("D9").Value = 1.025074385140119
$ws.Range("E9").Value = 1.044663009973872
$ws.Range("F9").Value = 1.04846843139572
$ws.Range("I9").Value = 1.028174127105039
$ws.Range("J9").Value = 1.0265078421839
$ws.Range("K9").Value = 1.028273020500442
$ws.Range("L9").Value = 1.047797584477989
$ws.Range("M9").Value = 1.051590864901751
$ws.Range("N9").Value = 1.012877047694932
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.018791090049523
$ws.Range("D10").Value = 1.023764413206515
$ws.Range("E10").Value = 1.042530645334649
$ws.Range("F10").Value = 1.046172427952731
$ws.Range("I10").Value = 1.027813759249521
$ws.Range("J10").Value = 1.025244000893299
$ws.Range("K10").Value = 1.027254305364646
$ws.Range("L10").Value = 1.04595307181217
$ws.Range("M10").Value = 1.04958207105437
$ws.Range("N10").Value = 1.012457137034687
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.018001711976763
$ws.Range("D11").Value = 1.023195427062574
$ws.Range("E11").Value = 1.041607205342582
$ws.Range("F11").Value = 1.045177991737064
$ws.Range("I11").Value = 1.02765475727549
$ws.Range("J11").Value = 1.024694630871763
$ws.Range("K11").Value = 1.026810657552251
$ws.Range("L11").Value = 1.045153645768094
$ws.Range("M11").Value = 1.048711401125997
$ws.Range("N11").Value = 1.012274549880541
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.017708261770759
$ws.Range("D12").Value = 1.022983814451011
$ws.Range("E12").Value = 1.041264176633615
$ws.Range("F12").Value = 1.044808570476918
$ws.Range("I12").Value = 1.027595251575528
$ws.Range("J12").Value = 1.024490249921284
$ws.Range("K12").Value = 1.026645484299389
$ws.Range("L12").Value = 1.044856587561958
$ws.Range("M12").Value = 1.048387863105283
$ws.Range("N12").Value = 1.012206613555048
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.017771218755331
$ws.Range("D13").Value = 1.02302921812728
$ws.Range("E13").Value = 1.041337758533698
$ws.Range("F13").Value = 1.044887814643798
$ws.Range("I13").Value = 1.027608035909091
$ws.Range("J13").Value = 1.024534104866177
$ws.Range("K13").Value = 1.026680931864942
$ws.Range("L13").Value = 1.044920312869221
$ws.Range("M13").Value = 1.048457269181277
$ws.Range("N13").Value = 1.012221191357342
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.017977460208448
$ws.Range("D14").Value = 1.023177940529335
$ws.Range("E14").Value = 1.04157885096074
$ws.Range("F14").Value = 1.045147456177603
$ws.Range("I14").Value = 1.02764984760264
$ws.Range("J14").Value = 1.024677743228811
$ws.Range("K14").Value = 1.026797012095308
$ws.Range("L14").Value = 1.045129093232763
$ws.Range("M14").Value = 1.048684660102507
$ws.Range("N14").Value = 1.012268936599304
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.018104500460798
$ws.Range("D15").Value = 1.02326953804145
$ws.Range("E15").Value = 1.041727392895119
$ws.Range("F15").Value = 1.045307423981248
$ws.Range("I15").Value = 1.027675550145868
$ws.Range("J15").Value = 1.024766201033931
$ws.Range("K15").Value = 1.026868482199895
$ws.Range("L15").Value = 1.045257714224955
$ws.Range("M15").Value = 1.048824745456939
$ws.Range("N15").Value = 1.012298338716129
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.018843445003604
$ws.Range("D16").Value = 1.02380213774553
$ws.Range("E16").Value = 1.042591928206425
$ws.Range("F16").Value = 1.046238419633622
$ws.Range("I16").Value = 1.0278242492713
$ws.Range("J16").Value = 1.025280415953034
$ws.Range("K16").Value = 1.02728369522218
$ws.Range("L16").Value = 1.046006111103548
$ws.Range("M16").Value = 1.049639836246331
$ws.Range("N16").Value = 1.012469238607592
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.019306541425591
$ws.Range("D17").Value = 1.024135751340795
$ws.Range("E17").Value = 1.04313419474542
$ws.Range("F17").Value = 1.046822337328495
$ws.Range("I17").Value = 1.027916731372271
$ws.Range("J17").Value = 1.025602400686275
$ws.Range("K17").Value = 1.027543467000066
$ws.Range("L17").Value = 1.04647535934997
$ws.Range("M17").Value = 1.050150890428715
$ws.Range("N17").Value = 1.012576234818937
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.019576506073797
$ws.Range("D18").Value = 1.024330172826834
$ws.Range("E18").Value = 1.043450479183671
$ws.Range("F18").Value = 1.047162902746759
$ws.Range("I18").Value = 1.027970389027276
$ws.Range("J18").Value = 1.025790004689417
$ws.Range("K18").Value = 1.027694742759831
$ws.Range("L18").Value = 1.0467489928386
$ws.Range("M18").Value = 1.05044889811016
$ws.Range("N18").Value = 1.012638570320924
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.019668531423645
$ws.Range("D19").Value = 1.02439643679176
$ws.Range("E19").Value = 1.043558322384466
$ws.Range("F19").Value = 1.047279022868382
$ws.Range("I19").Value = 1.02798863647564
$ws.Range("J19").Value = 1.025853938228675
$ws.Range("K19").Value = 1.027746282419273
$ws.Range("L19").Value = 1.046842282870489
$ws.Range("M19").Value = 1.050550497338915
$ws.Range("N19").Value = 1.012659812656559
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.019256871275835
$ws.Range("D20").Value = 1.024099975321516
$ws.Range("E20").Value = 1.043076015798347
$ws.Range("F20").Value = 1.046759690999401
$ws.Range("I20").Value = 1.027906838462252
$ws.Range("J20").Value = 1.025567875914424
$ws.Range("K20").Value = 1.027515621271149
$ws.Range("L20").Value = 1.046425020840372
$ws.Range("M20").Value = 1.050096067626286
$ws.Range("N20").Value = 1.012564762753836
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.017916733895694
$ws.Range("D21").Value = 1.023134152838604
$ws.Range("E21").Value = 1.041507855898271
$ws.Range("F21").Value = 1.045070999459516
$ws.Range("I21").Value = 1.027637547397492
$ws.Range("J21").Value = 1.024635454175259
$ws.Range("K21").Value = 1.026762839937158
$ws.Range("L21").Value = 1.045067615852692
$ws.Range("M21").Value = 1.048617702839328
$ws.Range("N21").Value = 1.012254880004688
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.017072745757513
$ws.Range("D22").Value = 1.022525363195278
$ws.Range("E22").Value = 1.040521761793005
$ws.Range("F22").Value = 1.044008998342289
$ws.Range("I22").Value = 1.027465656846115
$ws.Range("J22").Value = 1.024047348408056
$ws.Range("K22").Value = 1.026287321487087
$ws.Range("L22").Value = 1.044213489296377
$ws.Range("M22").Value = 1.047687427502906
$ws.Range("N22").Value = 1.012059376830869
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.017520292832154
$ws.Range("D23").Value = 1.022848240450188
$ws.Range("E23").Value = 1.041044522959676
$ws.Range("F23").Value = 1.044572011148755
$ws.Range("I23").Value = 1.027557023691446
$ws.Range("J23").Value = 1.024359290996432
$ws.Range("K23").Value = 1.026539613267675
$ws.Range("L23").Value = 1.044666343207098
$ws.Range("M23").Value = 1.048180658569627
$ws.Range("N23").Value = 1.012163080273146
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.019279315546958
$ws.Range("D24").Value = 1.024116141488991
$ws.Range("E24").Value = 1.043102304392276
$ws.Range("F24").Value = 1.04678799825245
$ws.Range("I24").Value = 1.027911309524903
$ws.Range("J24").Value = 1.02558347680375
$ws.Range("K24").Value = 1.027528204313596
$ws.Range("L24").Value = 1.046447766865409
$ws.Range("M24").Value = 1.050120839942668
$ws.Range("N24").Value = 1.012569946713579
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.021314626231961
$ws.Range("D25").Value = 1.025580818630742
$ws.Range("E25").Value = 1.045489852816928
$ws.Range("F25").Value = 1.049358610013046
$ws.Range("I25").Value = 1.028311237048313
$ws.Range("J25").Value = 1.026996071625081
$ws.Range("K25").Value = 1.028665814074539
$ws.Range("L25").Value = 1.048512237822902
$ws.Range("M25").Value = 1.052369129352266
$ws.Range("N25").Value = 1.013039208562512
